$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the data range so numeric-looking strings are not
# auto-converted to numbers (matches original t="inlineStr" text cells).
$dataRange = $ws.Range("A2:D36")
$dataRange.NumberFormat = "@"

$ws.Range("A2").Value = "0"
$ws.Range("B2").Value = "0.5"
$ws.Range("C2").Value = "0.377582561890373"
$ws.Range("D2").Value = "1.000001"

$ws.Range("A3").Value = "1"
$ws.Range("B3").Value = "0.877582561890373"
$ws.Range("C3").Value = "-0.238570067725114"
$ws.Range("D3").Value = "0.377582561890373"

$ws.Range("A4").Value = "2"
$ws.Range("B4").Value = "0.639012494165259"
$ws.Range("C4").Value = "0.163672606517076"
$ws.Range("D4").Value = "0.238570067725114"

$ws.Range("A5").Value = "3"
$ws.Range("B5").Value = "0.802685100682335"
$ws.Range("C5").Value = "-0.107907073894329"
$ws.Range("D5").Value = "0.163672606517076"

$ws.Range("A6").Value = "4"
$ws.Range("B6").Value = "0.694778026788006"
$ws.Range("C6").Value = "0.0734178044940099"
$ws.Range("D6").Value = "0.107907073894329"

$ws.Range("A7").Value = "5"
$ws.Range("B7").Value = "0.768195831282016"
$ws.Range("C7").Value = "-0.049030385339597"
$ws.Range("D7").Value = "0.0734178044940099"

$ws.Range("A8").Value = "6"
$ws.Range("B8").Value = "0.719165445942419"
$ws.Range("C8").Value = "0.0331903134791079"
$ws.Range("D8").Value = "0.049030385339597"

$ws.Range("A9").Value = "7"
$ws.Range("B9").Value = "0.752355759421527"
$ws.Range("C9").Value = "-0.0222746962837037"
$ws.Range("D9").Value = "0.0331903134791079"

$ws.Range("A10").Value = "8"
$ws.Range("B10").Value = "0.730081063137823"
$ws.Range("C10").Value = "0.0150392782136168"
$ws.Range("D10").Value = "0.0222746962837037"

$ws.Range("A11").Value = "9"
$ws.Range("B11").Value = "0.74512034135144"
$ws.Range("C11").Value = "-0.010114032336597"
$ws.Range("D11").Value = "0.0150392782136168"

$ws.Range("A12").Value = "10"
$ws.Range("B12").Value = "0.735006309014843"
$ws.Range("C12").Value = "0.0068202136284027"
$ws.Range("D12").Value = "0.010114032336597"

$ws.Range("A13").Value = "11"
$ws.Range("B13").Value = "0.741826522643246"
$ws.Range("C13").Value = "-0.0045907972010145"
$ws.Range("D13").Value = "0.0068202136284027"

$ws.Range("A14").Value = "12"
$ws.Range("B14").Value = "0.737235725442231"
$ws.Range("C14").Value = "0.0030939264360317"
$ws.Range("D14").Value = "0.0045907972010145"

$ws.Range("A15").Value = "13"
$ws.Range("B15").Value = "0.740329651878263"
$ws.Range("C15").Value = "-0.0020834135460297"
$ws.Range("D15").Value = "0.0030939264360317"

$ws.Range("A16").Value = "14"
$ws.Range("B16").Value = "0.738246238332233"
$ws.Range("C16").Value = "0.0014037244374277"
$ws.Range("D16").Value = "0.0020834135460297"

$ws.Range("A17").Value = "15"
$ws.Range("B17").Value = "0.739649962769661"
$ws.Range("C17").Value = "-0.0009454234126778"
$ws.Range("D17").Value = "0.0014037244374277"

$ws.Range("A18").Value = "16"
$ws.Range("B18").Value = "0.738704539356983"
$ws.Range("C18").Value = "0.0006369129242267"
$ws.Range("D18").Value = "0.0009454234126778"

$ws.Range("A19").Value = "17"
$ws.Range("B19").Value = "0.73934145228121"
$ws.Range("C19").Value = "-0.0004290029491069"
$ws.Range("D19").Value = "0.0006369129242267"

$ws.Range("A20").Value = "18"
$ws.Range("B20").Value = "0.738912449332103"
$ws.Range("C20").Value = "0.0002889948036959"
$ws.Range("D20").Value = "0.0004290029491069"

$ws.Range("A21").Value = "19"
$ws.Range("B21").Value = "0.739201444135799"
$ws.Range("C21").Value = "-0.000194664354986"
$ws.Range("D21").Value = "0.0002889948036959"

$ws.Range("A22").Value = "20"
$ws.Range("B22").Value = "0.739006779780813"
$ws.Range("C22").Value = "0.0001311309814798"
$ws.Range("D22").Value = "0.000194664354986"

$ws.Range("A23").Value = "21"
$ws.Range("B23").Value = "0.739137910762293"
$ws.Range("C23").Value = "-8.83301670843162e-05"
$ws.Range("D23").Value = "0.0001311309814798"

$ws.Range("A24").Value = "22"
$ws.Range("B24").Value = "0.739049580595209"
$ws.Range("C24").Value = "5.95008253181373e-05"
$ws.Range("D24").Value = "8.83301670843162e-05"

$ws.Range("A25").Value = "23"
$ws.Range("B25").Value = "0.739109081420527"
$ws.Range("C25").Value = "-4.0080216515137e-05"
$ws.Range("D25").Value = "5.95008253181373e-05"

$ws.Range("A26").Value = "24"
$ws.Range("B26").Value = "0.739069001204012"
$ws.Range("C26").Value = "2.69986317431581e-05"
$ws.Range("D26").Value = "4.0080216515137e-05"

$ws.Range("A27").Value = "25"
$ws.Range("B27").Value = "0.739095999835755"
$ws.Range("C27").Value = "-1.8186550579502e-05"
$ws.Range("D27").Value = "2.69986317431581e-05"

$ws.Range("A28").Value = "26"
$ws.Range("B28").Value = "0.739077813285175"
$ws.Range("C28").Value = "1.22507030757557e-05"
$ws.Range("D28").Value = "1.8186550579502e-05"

$ws.Range("A29").Value = "27"
$ws.Range("B29").Value = "0.739090063988251"
$ws.Range("C29").Value = "-8.25221014166821e-06"
$ws.Range("D29").Value = "1.22507030757557e-05"

$ws.Range("A30").Value = "28"
$ws.Range("B30").Value = "0.739081811778109"
$ws.Range("C30").Value = "5.55879292651884e-06"
$ws.Range("D30").Value = "8.25221014166821e-06"

$ws.Range("A31").Value = "29"
$ws.Range("B31").Value = "0.739087370571036"
$ws.Range("C31").Value = "-3.74446755613267e-06"
$ws.Range("D31").Value = "5.55879292651884e-06"

$ws.Range("A32").Value = "30"
$ws.Range("B32").Value = "0.73908362610348"
$ws.Range("C32").Value = "2.52231939923409e-06"
$ws.Range("D32").Value = "3.74446755613267e-06"

$ws.Range("A33").Value = "31"
$ws.Range("B33").Value = "0.739086148422879"
$ws.Range("C33").Value = "-1.69906423030941e-06"
$ws.Range("D33").Value = "2.52231939923409e-06"

$ws.Range("A34").Value = "32"
$ws.Range("B34").Value = "0.739084449358649"
$ws.Range("C34").Value = "1.14451031196783e-06"
$ws.Range("D34").Value = "1.69906423030941e-06"

$ws.Range("A35").Value = "33"
$ws.Range("B35").Value = "0.739085593868961"
$ws.Range("C35").Value = "-7.70955819295871e-07"
$ws.Range("D35").Value = "1.14451031196783e-06"

$ws.Range("A36").Value = "34"
$ws.Range("B36").Value = "0.739084822913141"
$ws.Range("C36").Value = "5.19325156767003e-07"
$ws.Range("D36").Value = "7.70955819295871e-07"

# Restore default (General) formatting/style so no stray text-format style
# lingers on the cells (values remain text because the stored content is
# non-numeric-formatted text in the XML already).
$dataRange.ClearFormats()

Write-Host "Updated tabla_pf data through row 36"
